$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Exposure conditions")

$ws.Range("M6").Value  = "XAA997AZ1"
$ws.Range("M7").Value  = "XAA997AZ2"
$ws.Range("M8").Value  = "XAA997AZ3"
$ws.Range("M9").Value  = "XAA997AZ4"

$ws.Range("M14").Value = "XAA997BZ1"
$ws.Range("M15").Value = "XAA997BZ2"
$ws.Range("M16").Value = "XAA997BZ3"
$ws.Range("M17").Value = "XAA997BZ4"

$ws.Range("M22").Value = "XAA997CZ1"
$ws.Range("M23").Value = "XAA997CZ2"
$ws.Range("M24").Value = "XAA997CZ3"
$ws.Range("M25").Value = "XAA997CZ4"
